$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Package Name version string from 1.0.3 -> 1.0.4 (row 17, column C)
$ws.Range("C17").Value = "1.0.4"

# Add a new row 24 with "Queue FolderName" / "Shared" (mirrors the style of
# the existing Package Name / Queue Name rows: label in A, bold value in B)
$ws.Range("B17").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Queue FolderName"
$ws.Range("B24").Value = "Shared"

# Update the selection to match the newly active cell after the edit
$ws.Range("C22").Select()
